# Remove the trailing "Ver no Jupiter ..." / "(c) 2020 ..." footer block
# (and the blank paragraph right before it) that followed the bibliography
# text at the end of the document.

$d = $word.ActiveDocument

$markerText = "Ver no Jupiter Salvar em pdf Salvar em docx"

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $markerText) {
        $target = $i
        break
    }
}

if ($null -eq $target) {
    throw "Could not find the 'Ver no Jupiter ...' paragraph"
}

# Paragraph right before the marker is the blank separator paragraph that
# should be removed together with the marker paragraph and the following
# copyright paragraph.
$firstToRemove = $target - 1
$lastToRemove = $target + 1

$startRange = $d.Paragraphs.Item($firstToRemove).Range.Start
$endRange = $d.Paragraphs.Item($lastToRemove).Range.End

$r = $d.Range($startRange, $endRange)
$r.Delete()
